# Weekly fruit/vegetable price update for "Pepino dulce" (Vega Modelo de Temuco).
# Two new weekly observations are inserted right after the existing row 248,
# pushing the previously-recorded rows 249-265 down to 251-267 (unchanged),
# while rows 247 and 248 are updated with the new observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 18   # columns A..R

# 1) Snapshot the current (pre-edit) contents of rows 247 and 248 - these
#    values get pushed down into the two newly inserted rows.
$row247 = @()
$row248 = @()
for ($c = 1; $c -le $lastCol; $c++) {
    $row247 += $ws.Cells.Item(247, $c).Value()
    $row248 += $ws.Cells.Item(248, $c).Value()
}

# 2) Insert two blank rows at 249:250, shifting old rows 249-265 down to 251-267.
$ws.Rows("249:250").Insert()

# 3) Re-populate the newly inserted rows 249 and 250 with the snapshot taken
#    from the old rows 247 and 248 (duplicated down, unedited).
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item(249, $c).Value = $row247[$c - 1]
    $ws.Cells.Item(250, $c).Value = $row248[$c - 1]
}

# 4) Apply the new-observation edits to row 247 (date, quality, volume).
$ws.Range("D247").Value = 44783
$ws.Range("I247").Value = "Extra"
$ws.Range("J247").Value = 50

# 5) Apply the new-observation edits to row 248 (date, volume).
$ws.Range("D248").Value = 44783
$ws.Range("J248").Value = 120
